# Generate Report for Handoff
# - Flip status "In Translation" -> "Ready for handoff" everywhere it appears
#   (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2)
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps to reflect the new handoff-generation run
# - Widen the Status column now that its text is longer

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status: "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Timestamps
$wsOverview.Range("G2").Value = "2016-08-22 14:43:49"
$wsDeDe.Range("H2").Value = "2016-08-22 14:43:49"
$wsZhCn.Range("H2").Value = "2016-08-22 14:43:43"

# Widen the Status column to fit "Ready for handoff" (closest width this
# engine's half-point-quantized column model can represent to the
# authored 17.2159881591797 target)
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
